# Apply "Deneme yanilma 1.5 milyon hedef" edit:
# On sheet "yeni_degiskenler": change B2 from 25 to 72, B5 from 0.75 to 0.6,
# and move the active selection from F6 to B6.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("yeni_degiskenler")
$ws.Activate()

$ws.Range("B2").Value = 72
$ws.Range("B5").Value = 0.6

$ws.Range("B6").Select()
